$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 11 - this shifts the existing SUM row (12) and
# average row (13) down to 13 and 14, and Excel auto-updates their formulas'
# cell references (SUM(B2:B11) -> SUM(B2:B12), B12 -> B13).
$ws.Rows("11").Insert()

# Populate the new row 11 with the new timesheet entry.
$ws.Range("A11").Value = 40882
$ws.Range("A11").NumberFormat = "d-mmm"
$ws.Range("A11").HorizontalAlignment = -4152
$ws.Range("B11").Value = 30
$ws.Range("C11").Value = "Replace images on Curriculum and Contact Us. Released  v.2011.12.05.1"

# Update the selection to match the post-edit cursor position.
$ws.Range("B15").Select()
